$d = $word.ActiveDocument

# Locate the target paragraph: the header line that currently contains the
# single-run "Portfolio: ... | GitHub: ... | LinkedIn: ..." text, wherever it
# happens to sit in the document.
$findRange = $d.Content
$found = $findRange.Find.Execute("Portfolio:*LinkedIn:*", $true, $false, $true,
                                  $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the Portfolio/GitHub/LinkedIn header paragraph"
}
$p = $findRange.Paragraphs(1)
$r = $p.Range

$frag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:after="120"/>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">Portfolio: </w:t>
            </w:r>
            <w:hyperlink r:id="rIdPortfolio">
              <w:r>
                <w:rPr>
                  <w:color w:val="0563C1"/>
                  <w:u w:val="single"/>
                </w:rPr>
                <w:t>slimeq.github.io</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">  |  GitHub: </w:t>
            </w:r>
            <w:hyperlink r:id="rIdGitHub">
              <w:r>
                <w:rPr>
                  <w:color w:val="0563C1"/>
                  <w:u w:val="single"/>
                </w:rPr>
                <w:t>github.com/SlimeQ</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve">  |  LinkedIn: </w:t>
            </w:r>
            <w:hyperlink r:id="rIdLinkedIn">
              <w:r>
                <w:rPr>
                  <w:color w:val="0563C1"/>
                  <w:u w:val="single"/>
                </w:rPr>
                <w:t>quincy-campbell-131559b2</w:t>
              </w:r>
            </w:hyperlink>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rIdPortfolio" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://slimeq.github.io/" TargetMode="External"/>
        <Relationship Id="rIdGitHub" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/SlimeQ" TargetMode="External"/>
        <Relationship Id="rIdLinkedIn" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://linkedin.com/in/quincy-campbell-131559b2" TargetMode="External"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($frag)
Write-Output $p.Range.Text
